$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Text = "29.458.51" },
    @{ Ref = "E2"; Text = "  -3.10%  " },
    @{ Ref = "D3"; Text = "1.992.48" },
    @{ Ref = "E3"; Text = "  -4.96%  " },
    @{ Ref = "E4"; Text = "  +1.28%  " },
    @{ Ref = "D5"; Text = "329.01" },
    @{ Ref = "E5"; Text = "  -4.07%  " },
    @{ Ref = "E6"; Text = "  +1.13%  " },
    @{ Ref = "E7"; Text = "  -4.64%  " },
    @{ Ref = "D8"; Text = "0.4228" },
    @{ Ref = "E8"; Text = "  -4.53%  " },
    @{ Ref = "D9"; Text = "54.15" },
    @{ Ref = "E9"; Text = "  -0.77%  " },
    @{ Ref = "D10"; Text = "0.08919" },
    @{ Ref = "E10"; Text = "  -4.72%  " },
    @{ Ref = "E11"; Text = "  -5.13%  " },
    @{ Ref = "D12"; Text = "23.22" },
    @{ Ref = "E12"; Text = "  -6.40%  " },
    @{ Ref = "D13"; Text = "1.992.00" },
    @{ Ref = "E13"; Text = "  -3.46%  " },
    @{ Ref = "D14"; Text = "7.970" },
    @{ Ref = "E14"; Text = "  -7.09%  " },
    @{ Ref = "D15"; Text = "6.456" },
    @{ Ref = "E15"; Text = "  -6.81%  " },
    @{ Ref = "D16"; Text = "1.015" },
    @{ Ref = "E16"; Text = "  +1.24%  " },
    @{ Ref = "D17"; Text = "94.12" },
    @{ Ref = "E17"; Text = "  -7.18%  " },
    @{ Ref = "D18"; Text = "0.00001110" },
    @{ Ref = "E18"; Text = "  -4.19%  " },
    @{ Ref = "D19"; Text = "0.06759" },
    @{ Ref = "E19"; Text = "  +1.21%  " },
    @{ Ref = "D20"; Text = "19.41" },
    @{ Ref = "E20"; Text = "  -8.43%  " },
    @{ Ref = "D21"; Text = "1.013" },
    @{ Ref = "E21"; Text = "  +1.20%  " },
    @{ Ref = "D22"; Text = "5.934" },
    @{ Ref = "E22"; Text = "  -6.27%  " },
    @{ Ref = "D23"; Text = "29.489.79" },
    @{ Ref = "E23"; Text = "  -3.10%  " },
    @{ Ref = "E24"; Text = "  -3.98%  " },
    @{ Ref = "D25"; Text = "2.329" },
    @{ Ref = "E25"; Text = "  +0.90%  " },
    @{ Ref = "D27"; Text = "156.87" },
    @{ Ref = "E27"; Text = "  -3.73%  " },
    @{ Ref = "D28"; Text = "6.254" },
    @{ Ref = "E28"; Text = "  -8.17%  " },
    @{ Ref = "D29"; Text = "2.302" },
    @{ Ref = "E29"; Text = "  -8.42%  " },
    @{ Ref = "D30"; Text = "127.58" },
    @{ Ref = "E30"; Text = "  -4.49%  " },
    @{ Ref = "D31"; Text = "1.057" },
    @{ Ref = "E31"; Text = "  -7.29%  " },
    @{ Ref = "D32"; Text = "0.09932" },
    @{ Ref = "E32"; Text = "  -5.39%  " },
    @{ Ref = "D33"; Text = "1.542" },
    @{ Ref = "E33"; Text = "  -6.90%  " },
    @{ Ref = "D34"; Text = "5.831" },
    @{ Ref = "E34"; Text = "  -7.09%  " },
    @{ Ref = "D35"; Text = "3.792" },
    @{ Ref = "E35"; Text = "  -2.18%  " },
    @{ Ref = "E36"; Text = "  -6.78%  " },
    @{ Ref = "D37"; Text = "9.206" },
    @{ Ref = "E37"; Text = "  -9.59%  " },
    @{ Ref = "D38"; Text = "0.06372" },
    @{ Ref = "E38"; Text = "  -6.48%  " },
    @{ Ref = "D39"; Text = "1.297" },
    @{ Ref = "E39"; Text = "  -3.36%  " },
    @{ Ref = "D40"; Text = "0.6536" },
    @{ Ref = "E40"; Text = "  -6.74%  " },
    @{ Ref = "E41"; Text = "  -7.98%  " },
    @{ Ref = "D42"; Text = "0.2038" },
    @{ Ref = "E42"; Text = "  -8.13%  " },
    @{ Ref = "D43"; Text = "1.013" },
    @{ Ref = "E43"; Text = "  +1.15%  " },
    @{ Ref = "D44"; Text = "0.6328" },
    @{ Ref = "E44"; Text = "  -7.68%  " },
    @{ Ref = "D45"; Text = "13.53" },
    @{ Ref = "E45"; Text = "  -5.88%  " },
    @{ Ref = "D46"; Text = "2.206" },
    @{ Ref = "E46"; Text = "  -5.95%  " },
    @{ Ref = "D47"; Text = "1.310" },
    @{ Ref = "E47"; Text = "  -5.39%  " },
    @{ Ref = "D48"; Text = "3.503" },
    @{ Ref = "E48"; Text = "  -3.59%  " },
    @{ Ref = "D49"; Text = "0.00000000342" },
    @{ Ref = "E49"; Text = "  -3.60%  " },
    @{ Ref = "D50"; Text = "0.06954" },
    @{ Ref = "E50"; Text = "  -4.14%  " },
    @{ Ref = "D51"; Text = "1.131" },
    @{ Ref = "E51"; Text = "  -7.98%  " }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Ref)
    $c.NumberFormat = "@"
    $c.Value = $u.Text
}

